$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 19 values (columns G, H, I, J go from 2 -> 5)
$ws.Range("G19").Value = 5
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = 5

# M19 goes from 4 -> 5
$ws.Range("M19").Value = 5

# L19 is a shared formula SUM(C19:J19); recalculated automatically to 40.

# Update the active selection on the sheet to O19
$ws.Range("O19").Select()
